# Updates the cryptos list (price / 1h volume columns, plus a ranking swap
# between dogwifhat and Hedera at rows 37-38) to match the latest scrape.
#
# Note: several Price values (column D) are strings that happen to look like
# plain decimals (e.g. "6.50", "45.60") and must stay TEXT so the trailing
# zero / exact formatting survives. A leading apostrophe forces Excel to
# store them as text (quote-prefixed) instead of silently coercing them to
# numbers, which would drop the trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.665.40'
$ws.Cells.Item(2, 5).Value = '  -1.10%  '

$ws.Cells.Item(3, 4).Value = '3.770.48'
$ws.Cells.Item(3, 5).Value = '  -1.99%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = '''596.13'
$ws.Cells.Item(5, 5).Value = '  -0.98%  '

$ws.Cells.Item(6, 4).Value = '''168.17'
$ws.Cells.Item(6, 5).Value = '  -0.49%  '

$ws.Cells.Item(7, 4).Value = '3.768.16'
$ws.Cells.Item(7, 5).Value = '  -2.02%  '

$ws.Cells.Item(8, 5).Value = '  +0.04%  '

$ws.Cells.Item(9, 5).Value = '  -1.00%  '

$ws.Cells.Item(10, 5).Value = '  -0.88%  '

$ws.Cells.Item(11, 4).Value = '''6.50'
$ws.Cells.Item(11, 5).Value = '  +0.41%  '

$ws.Cells.Item(12, 4).Value = '''0.451'

$ws.Cells.Item(13, 5).Value = '  +4.06%  '

$ws.Cells.Item(14, 4).Value = '''36.17'
$ws.Cells.Item(14, 5).Value = '  -2.48%  '

$ws.Cells.Item(15, 4).Value = '4.404.78'
$ws.Cells.Item(15, 5).Value = '  -2.01%  '

$ws.Cells.Item(16, 4).Value = '3.768.28'
$ws.Cells.Item(16, 5).Value = '  -1.98%  '

$ws.Cells.Item(17, 4).Value = '''18.48'
$ws.Cells.Item(17, 5).Value = '  -0.22%  '

$ws.Cells.Item(18, 4).Value = '67.572.71'
$ws.Cells.Item(18, 5).Value = '  -1.38%  '

$ws.Cells.Item(19, 5).Value = '  -2.80%  '

$ws.Cells.Item(21, 4).Value = '''10.45'
$ws.Cells.Item(21, 5).Value = '  -6.55%  '

$ws.Cells.Item(22, 4).Value = '''465.98'
$ws.Cells.Item(22, 5).Value = '  -1.11%  '

$ws.Cells.Item(23, 4).Value = '''0.715'
$ws.Cells.Item(23, 5).Value = '  -2.48%  '

$ws.Cells.Item(24, 4).Value = '''0.0000147'
$ws.Cells.Item(24, 5).Value = '  -7.99%  '

$ws.Cells.Item(25, 4).Value = '''83.56'
$ws.Cells.Item(25, 5).Value = '  +0.05%  '

$ws.Cells.Item(26, 5).Value = '  -2.05%  '

$ws.Cells.Item(27, 4).Value = '''12.08'
$ws.Cells.Item(27, 5).Value = '  -0.22%  '

$ws.Cells.Item(28, 4).Value = '''10.27'
$ws.Cells.Item(28, 5).Value = '  +0.38%  '

$ws.Cells.Item(29, 5).Value = '  -0.04%  '

$ws.Cells.Item(30, 5).Value = '  -1.75%  '

$ws.Cells.Item(31, 4).Value = '3.922.97'
$ws.Cells.Item(31, 5).Value = '  -1.87%  '

$ws.Cells.Item(32, 4).Value = '''7.58'
$ws.Cells.Item(32, 5).Value = '  -1.67%  '

$ws.Cells.Item(33, 4).Value = '''30.45'

$ws.Cells.Item(34, 5).Value = '  -3.82%  '

$ws.Cells.Item(35, 4).Value = '''9.09'

$ws.Cells.Item(36, 4).Value = '3.734.99'
$ws.Cells.Item(36, 5).Value = '  -2.08%  '

$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(37, 4).Value = '''0.103'
$ws.Cells.Item(37, 5).Value = '  -1.16%  '

$ws.Cells.Item(38, 2).Value = 'dogwifhat'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(38, 4).Value = '''3.68'
$ws.Cells.Item(38, 5).Value = '  -2.45%  '

$ws.Cells.Item(39, 5).Value = '  -1.79%  '

$ws.Cells.Item(40, 5).Value = '  -1.44%  '

$ws.Cells.Item(41, 4).Value = '''5.78'
$ws.Cells.Item(41, 5).Value = '  -2.73%  '

$ws.Cells.Item(43, 5).Value = '  -1.68%  '

$ws.Cells.Item(45, 4).Value = '''8.63'
$ws.Cells.Item(45, 5).Value = '  -0.99%  '

$ws.Cells.Item(46, 5).Value = '  -2.63%  '

$ws.Cells.Item(47, 4).Value = '''45.60'
$ws.Cells.Item(47, 5).Value = '  -2.98%  '

$ws.Cells.Item(48, 4).Value = '''394.94'
$ws.Cells.Item(48, 5).Value = '  -5.65%  '

$ws.Cells.Item(49, 4).Value = '''0.000267'
$ws.Cells.Item(49, 5).Value = '  -9.03%  '

$ws.Cells.Item(50, 5).Value = '  -1.16%  '

$ws.Cells.Item(51, 4).Value = '''39.12'
$ws.Cells.Item(51, 5).Value = '  +2.72%  '
